# Add a new "UK" Test Data sheet, cloned from the existing "Poland" sheet,
# then adjust its content to describe the UK market test data.

$wb = $excel.ActiveWorkbook

# Poland is the template for every "market" sheet in this workbook.
$template = $wb.Worksheets.Item("Poland")

# Clone it and place the copy right after "Poland" (i.e. as the new last sheet).
$template.Copy($null, $template)
$newWs = $wb.ActiveSheet
$newWs.Name = "UK"

# Update the market-specific identifiers (NGC code first, then the market
# name, so new shared-string entries are appended in that order).
$newWs.Range("B4").Value = "NGC-2741/T3345/T3343/T3342"
$newWs.Range("B2").Value = "UK Market"

# Insert a new row for the "GMPIM" entry above the "Wg"/"Miscellaneous" rows,
# copying the formatting of the row above it so the new cell keeps the same
# border style as its neighbours.
$newWs.Rows("9:9").Insert()
$newWs.Range("A8").Copy()
$newWs.Range("A9").PasteSpecial(-4122)
$newWs.Range("A9").Value = "GMPIM"

# Make the new sheet the active one, with A9 selected, mirroring how it was
# left after the edit was made in Excel.
[void]$newWs.Range("A9").Select()
